$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 129, pushing existing rows 129..154 down to 130..155
$ws.Rows.Item(129).Insert()

# Populate the new row 129 with the new weekly record
$ws.Cells.Item(129, 1).Value = 8
$ws.Cells.Item(129, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(129, 3).Value = "Coquimbo"
$ws.Cells.Item(129, 4).Value = 44782
$ws.Cells.Item(129, 5).Value = 4
$ws.Cells.Item(129, 6).Value = 100112001
$ws.Cells.Item(129, 7).Value = "Berenjena"
$ws.Cells.Item(129, 8).Value = "Sin especificar"
$ws.Cells.Item(129, 9).Value = "Primera"
$ws.Cells.Item(129, 10).Value = 540
$ws.Cells.Item(129, 11).Value = 10000
$ws.Cells.Item(129, 12).Value = 11000
$ws.Cells.Item(129, 13).Value = 10500
$ws.Cells.Item(129, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(129, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(129, 16).Value = 210
$ws.Cells.Item(129, 17).Value = 50
$ws.Cells.Item(129, 18).Value = "Hortaliza"
